$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A (element) and Column B (type) values, rows 2-10
$ws.Range("A2").Value = "demand1"
$ws.Range("B2").Value = "demand"

$ws.Range("A3").Value = "net1"
$ws.Range("B3").Value = "net"

$ws.Range("A4").Value = "pv1"
$ws.Range("B4").Value = "pv"

$ws.Range("A5").Value = "bat1"
$ws.Range("B5").Value = "bat"

$ws.Range("A6").Value = "CHP1"
$ws.Range("B6").Value = "CHP"

$ws.Range("A7").Value = "solar_th1"
$ws.Range("B7").Value = "solar_th"

$ws.Range("A8").Value = "pvt1"
$ws.Range("B8").Value = "pvt"

$ws.Range("A9").Value = "charging_station1"
$ws.Range("B9").Value = "charging_station"

$ws.Range("A10").Value = "charging_station2"
$ws.Range("B10").Value = "charging_station"
